$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data in column C (felelős) and D (megjegyzés) for the Feltöltés /
# Letöltés rows (7 and 8): "Niki" is responsible, with a note attached to
# the (to-be-merged) D7:D8 block.
$ws.Range("C7").Value = "Niki"
$ws.Range("C8").Value = "Niki"
$ws.Range("D7").Value = "esetleg online tárhellyel való megoldás keresése"

# Format the note cell: centered, wrapped text (new cellXfs style).
$d7 = $ws.Range("D7")
$d7.HorizontalAlignment = -4108   # xlCenter
$d7.WrapText = $true

# Mirror the exact same formatting onto D8 (copy formats only, no value)
# so both cells share one style before/while merging.
$d7.Copy()
$ws.Range("D8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Merge the note cell across rows 7-8.
$ws.Range("D7:D8").Merge()

# Widen column D to fit the longer note text.
$ws.Columns.Item(4).ColumnWidth = 24.6

# Update the active selection to the merged note cell.
$ws.Range("D7:D8").Select() | Out-Null
